# Apply updated cryptocurrency price/volume data to the worksheet
# (mirrors the scheduled GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.310.40'
$ws.Range("E2").Value = '  +2.62%  '
$ws.Range("D3").Value = '3.122.38'
$ws.Range("E3").Value = '  +1.03%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.17'
$ws.Range("E5").Value = '  +2.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '622.41'
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.04'
$ws.Range("E7").Value = '  +28.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.374'
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("D10").Value = '3.116.81'
$ws.Range("E10").Value = '  +0.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.754'
$ws.Range("E11").Value = '  +22.65%  '
$ws.Range("E12").Value = '  +6.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000252'
$ws.Range("E13").Value = '  +4.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.86'
$ws.Range("E14").Value = '  +8.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.55'
$ws.Range("E15").Value = '  +4.78%  '
$ws.Range("D16").Value = '91.170.96'
$ws.Range("E16").Value = '  +2.81%  '
$ws.Range("D17").Value = '3.693.16'
$ws.Range("E17").Value = '  +0.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.88'
$ws.Range("E18").Value = '  +14.60%  '
$ws.Range("D19").Value = '3.110.73'
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000216'
$ws.Range("E20").Value = '  +2.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.13'
$ws.Range("E21").Value = '  +5.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '442.92'
$ws.Range("E22").Value = '  +4.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.85'
$ws.Range("E23").Value = '  +7.02%  '
$ws.Range("E24").Value = '  +5.66%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.25'
$ws.Range("E25").Value = '  +10.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '89.15'
$ws.Range("E26").Value = '  +8.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.37'
$ws.Range("E27").Value = '  +3.98%  '
$ws.Range("D28").Value = '3.281.92'
$ws.Range("E28").Value = '  +1.44%  '
$ws.Range("E29").Value = '  -0.62%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.166'
$ws.Range("E30").Value = '  -2.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.23'
$ws.Range("E31").Value = '  +13.76%  '
$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.974'
$ws.Range("E32").Value = '  -9.91%  '
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '524.53'
$ws.Range("E33").Value = '  +2.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.73'
$ws.Range("E34").Value = '  +1.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.148'
$ws.Range("E35").Value = '  +13.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.08'
$ws.Range("E36").Value = '  +4.96%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '24.14'
$ws.Range("E37").Value = '  +8.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.30'
$ws.Range("E38").Value = '  +3.73%  '
$ws.Range("E39").Value = '  +3.57%  '
$ws.Range("B40").Value = 'WhiteBITCoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.28'
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0859'
$ws.Range("E41").Value = '  +23.99%  '
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.159'
$ws.Range("E42").Value = '  +20.89%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.397'
$ws.Range("E44").Value = '  +9.28%  '
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.94'
$ws.Range("E46").Value = '  +5.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '149.07'
$ws.Range("E47").Value = '  +2.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '44.16'
$ws.Range("E48").Value = '  +2.12%  '
$ws.Range("E49").Value = '  +7.68%  '
$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.28'
$ws.Range("E50").Value = '  +8.48%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '168.31'
$ws.Range("E51").Value = '  +4.11%  '
